$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "27.070.29"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.675.88"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue $ws.Range("D5") "215.39"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +2.31%  "
Set-TextValue $ws.Range("D9") "21.27"
$ws.Range("E9").Value = "  +5.01%  "
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "1.911.38"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "1.675.48"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("E14").Value = "  +1.00%  "
Set-TextValue $ws.Range("D15") "0.536"
$ws.Range("E15").Value = "  +1.79%  "
Set-TextValue $ws.Range("D16") "66.02"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "27.053.37"
$ws.Range("E17").Value = "  +0.52%  "
Set-TextValue $ws.Range("D18") "237.19"
$ws.Range("E18").Value = "  +1.80%  "
Set-TextValue $ws.Range("D19") "8.15"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").Value = "0.0₃0739"
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("E24").Value = "  -1.89%  "
Set-TextValue $ws.Range("D25") "147.12"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("E26").Value = "  +1.49%  "
Set-TextValue $ws.Range("D27") "16.35"
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").Value = "1.531.64"
$ws.Range("E33").Value = "  +5.07%  "
Set-TextValue $ws.Range("D34") "3.18"
$ws.Range("E34").Value = "  +1.69%  "
$ws.Range("E35").Value = "  +3.75%  "
Set-TextValue $ws.Range("D36") "2.39"
$ws.Range("E36").Value = "  -0.90%  "
Set-TextValue $ws.Range("D37") "0.594"
$ws.Range("E37").Value = "  +2.17%  "
Set-TextValue $ws.Range("D38") "0.918"
$ws.Range("E38").Value = "  +2.17%  "
$ws.Range("E39").Value = "  +2.32%  "
$ws.Range("E40").Value = "  +2.69%  "
$ws.Range("E41").Value = "  +0.07%  "
Set-TextValue $ws.Range("D42") "67.57"
$ws.Range("E42").Value = "  +2.01%  "
Set-TextValue $ws.Range("D43") "5.51"
$ws.Range("E43").Value = "  -3.83%  "
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").Value = "1.819.09"
$ws.Range("E46").Value = "  +0.27%  "
Set-TextValue $ws.Range("D47") "90.75"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("E49").Value = "  +2.32%  "
Set-TextValue $ws.Range("D50") "8.02"
$ws.Range("E50").Value = "  +5.37%  "
$ws.Range("E51").Value = "  +0.63%  "
